# Fruta / hortaliza, semanal
#
# A new weekly record is inserted at row 46 (Macroferia Regional de Talca -
# Mango), pushing all the subsequent records down by one row. The new
# record reuses the same market/product/category metadata as the record
# that used to be first in that block, but carries its own date, volume,
# price and $/Kg figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 46, shifting rows 46:101
# down to 47:102 (values, formatting and all).
$ws.Rows("46:46").Insert()

# Populate the newly inserted row with the new weekly entry. Columns that
# are not price/volume/date figures keep the same values the block already
# used (market, region, product taxonomy, unit, origin, kg/unit).
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 44580
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100108
$ws.Range("H46").Value = "Tropicales y subtropicales"
$ws.Range("I46").Value = 100108002
$ws.Range("J46").Value = "Mango"
$ws.Range("K46").Value = "Sin especificar"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 320
$ws.Range("N46").Value = 7000
$ws.Range("O46").Value = 7000
$ws.Range("P46").Value = 7000
$ws.Range("Q46").Value = "$/bandeja 4 kilos"
$ws.Range("R46").Value = "Perú"
$ws.Range("S46").Value = 1750
$ws.Range("T46").Value = 4
